# Apply the "Add files via upload" edit: populate the new boolean "aero"
# flag column (E) for every data row (2-38) on Sheet1, and update the
# sheet's view/selection state to match the resaved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Boolean values for column E ("aero"), keyed by row number, taken from
# the target workbook.
$aeroValues = @{
    2  = $true
    3  = $true
    4  = $false
    5  = $false
    6  = $false
    7  = $true
    8  = $true
    9  = $true
    10 = $true
    11 = $true
    12 = $false
    13 = $false
    14 = $true
    15 = $false
    16 = $false
    17 = $false
    18 = $false
    19 = $false
    20 = $false
    21 = $false
    22 = $false
    23 = $false
    24 = $false
    25 = $true
    26 = $false
    27 = $false
    28 = $false
    29 = $false
    30 = $false
    31 = $false
    32 = $false
    33 = $false
    34 = $false
    35 = $false
    36 = $false
    37 = $true
    38 = $false
}

foreach ($row in $aeroValues.Keys | Sort-Object) {
    $ws.Range("E$row").Value = $aeroValues[$row]
}

# Update the view so the selection / scroll position reflect the resave.
$ws.Range("E38").Select()
